$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws2 = $wb.Worksheets.Item("Delete Devices Loop A")

# ---------------------------------------------------------------
# Sheet 1: "Add Devices Loop A"
# Add two new "Volt Drop" columns (E for header rows 3/4, H/I for the
# per-device table) driven by new "Loading Details" methods.
# ---------------------------------------------------------------

# New table header cells H5/I5 (same style as the other headers in row 5)
# NOTE: these are entered first so "Loading Details Name..." is added to
# the shared string table ahead of the "Volt Drop..." strings below.
$ws1.Range("A5").Copy()
$ws1.Range("H5").PasteSpecial(-4122)
$ws1.Range("H5").Value = "Loading Details Name 1"

$ws1.Range("A5").Copy()
$ws1.Range("I5").PasteSpecial(-4122)
$ws1.Range("I5").Value = "Loading Details Name2"

# New merged-looking header cells E3/E4 (same look as the G column boxes)
$ws1.Range("G3").Copy()
$ws1.Range("E3").PasteSpecial(-4122)
$ws1.Range("E3").Value = "Volt Drop (V)"
$ws1.Range("E3").HorizontalAlignment = -4131
$ws1.Range("E3").WrapText = $true

$ws1.Range("G3").Copy()
$ws1.Range("E4").PasteSpecial(-4122)
$ws1.Range("E4").Value = "Volt Drop (worst case)"
$ws1.Range("E4").HorizontalAlignment = -4131
$ws1.Range("E4").WrapText = $true

# New data cells H6:I8 (same style as the other shaded boxes, e.g. G3)
$rows = 6,7,8
foreach ($r in $rows) {
    $ws1.Range("G3").Copy()
    $ws1.Range("H$r").PasteSpecial(-4122)
    $ws1.Range("H$r").Value = "Volt Drop (V)"
    $ws1.Range("H$r").HorizontalAlignment = -4131
    $ws1.Range("H$r").WrapText = $true

    $ws1.Range("G3").Copy()
    $ws1.Range("I$r").PasteSpecial(-4122)
    $ws1.Range("I$r").Value = "Volt Drop (worst case)"
    $ws1.Range("I$r").HorizontalAlignment = -4131
    $ws1.Range("I$r").WrapText = $true
}

# Rows 3 & 4 grow taller to fit the wrapped "Volt Drop" labels
$ws1.Rows.Item(3).RowHeight = 28.8
$ws1.Rows.Item(4).RowHeight = 43.2

# Column H widens to fit the new "Loading Details Name 1" text
$ws1.Range("H5").EntireColumn.AutoFit()

# ---------------------------------------------------------------
# Sheet 2: "Delete Devices Loop A"
# Rename the "Voltage Drop" labels to the shorter "Volt Drop" wording
# ---------------------------------------------------------------
$ws2.Range("A7").Value = "Volt Drop (V)"
$ws2.Range("C7").Value = "Volt Drop (worst case)"

# ---------------------------------------------------------------
# View state: "Delete Devices Loop A" becomes the active/selected tab,
# while "Add Devices Loop A" keeps a selection on the new E3:E4 cells.
# ---------------------------------------------------------------
$ws1.Range("E3:E4").Select()
$ws2.Activate()
$ws2.Range("C7").Select()
